$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Niveau espece")
$ws.Activate()
$ws.Rows.Item(10).Delete()
$ws.Range("A10:XFD10").Select()
